# Commit: "trying to have common splash screen for login"
#
# This adds a new employee record (with a duplicate "saved twice" row,
# presumably from the splash/login flow re-submitting) to EMP_RECORD,
# and logs matching login/attendance timestamps for those employee IDs
# on the EMP_ATTENDANCE sheet.

$wb = $excel.ActiveWorkbook

# --- EMP_RECORD: populate A1:I3 with the new employee rows ---
$empRecord = $wb.Worksheets.Item("EMP_RECORD")

$empRows = @(
    @("1", "32",  "Chetan",        "Chinchulkar", "safdsf", "asdf", "fsadf", "dfdsa", "2"),
    @("1", "32",  "Chetan",        "Chinchulkar", "safdsf", "asdf", "fsadf", "dfdsa", "2"),
    @("2", "323", "Chetansdafdsf", "Chinchulkar", "safdsf", "asdf", "fsadf", "dfdsa", "2")
)

for ($r = 0; $r -lt $empRows.Length; $r++) {
    $rowValues = $empRows[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $cell = $empRecord.Cells.Item($r + 1, $c + 1)
        # Force text storage so values like "1" / "32" stay strings, not numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $rowValues[$c]
    }
}

# --- EMP_ATTENDANCE: append login/attendance rows for the new employees ---
$empAttendance = $wb.Worksheets.Item("EMP_ATTENDANCE")

$attendanceRows = @(
    @("32",  "03-07-2022", "21:27:22"),
    @("323", "03-07-2022", "21:28:04")
)

$startRow = 3
for ($r = 0; $r -lt $attendanceRows.Length; $r++) {
    $rowValues = $attendanceRows[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $cell = $empAttendance.Cells.Item($startRow + $r, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $rowValues[$c]
    }
}
